# Quarterly indexing bug-fix: the dates in column A (rows 2-73) were
# recorded as the 1st of the reference month instead of the 15th of the
# *following* month. Shift every date forward by one month and pin the
# day-of-month to the 15th.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 73; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $old = $cell.Value2
    if ($old -eq $null) { continue }

    $d = [datetime]::FromOADate($old)
    $d = $d.AddMonths(1)

    $newSerial = $excel.Evaluate("DATE($($d.Year),$($d.Month),15)")
    $cell.Value2 = $newSerial
}
